$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 2.82
$ws.Range("L2").Value = 1.3
$ws.Range("U2").Value = 2.28
$ws.Range("V2").Value = 1.55
$ws.Range("AF2").Value = 24
$ws.Range("AK2").Value = 32

# Row 3
$ws.Range("F3").Value = 1.9

# Row 4
$ws.Range("G4").Value = 2.52
$ws.Range("I4").Value = 3.25
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 4.5
$ws.Range("O4").Value = 1.23
$ws.Range("R4").Value = 1.47
$ws.Range("S4").Value = 2.72
$ws.Range("T4").Value = 1.6
$ws.Range("U4").Value = 2.32
$ws.Range("V4").Value = 1.45
$ws.Range("W4").Value = 1.66
$ws.Range("X4").Value = 24
$ws.Range("Y4").Value = 18.5
$ws.Range("Z4").Value = 28
$ws.Range("AA4").Value = 60
$ws.Range("AB4").Value = 16
$ws.Range("AC4").Value = 11
$ws.Range("AD4").Value = 16
$ws.Range("AE4").Value = 38
$ws.Range("AF4").Value = 21
$ws.Range("AG4").Value = 14.5
$ws.Range("AH4").Value = 19
$ws.Range("AI4").Value = 46
$ws.Range("AJ4").Value = 40
$ws.Range("AK4").Value = 29
$ws.Range("AL4").Value = 40
$ws.Range("AM4").Value = 85
$ws.Range("AN4").Value = 18.5
$ws.Range("AO4").Value = 29

# Row 5
$ws.Range("S5").Value = 2.18
$ws.Range("T5").Value = 1.54
$ws.Range("X5").Value = 32

# Row 6
$ws.Range("L6").Value = 1.27
$ws.Range("M6").Value = 1.06
$ws.Range("R6").Value = 1.38
$ws.Range("S6").Value = 3.05
$ws.Range("T6").Value = 1.66
$ws.Range("U6").Value = 2.2
$ws.Range("V6").Value = 1.56
$ws.Range("X6").Value = 19.5
$ws.Range("Y6").Value = 12.5
$ws.Range("Z6").Value = 19
$ws.Range("AA6").Value = 40
$ws.Range("AB6").Value = 13.5
$ws.Range("AC6").Value = 8.800000000000001
$ws.Range("AD6").Value = 13
$ws.Range("AE6").Value = 29
$ws.Range("AF6").Value = 21
$ws.Range("AG6").Value = 13.5
$ws.Range("AH6").Value = 17
$ws.Range("AI6").Value = 40
$ws.Range("AJ6").Value = 46
$ws.Range("AK6").Value = 32
$ws.Range("AL6").Value = 42
$ws.Range("AM6").Value = 85
$ws.Range("AN6").Value = 26
$ws.Range("AO6").Value = 23

# Row 7
$ws.Range("H7").Value = 1.01
$ws.Range("I7").Value = 7.4
$ws.Range("J7").Value = 2.88
$ws.Range("K7").Value = 5.1
$ws.Range("L7").Value = 1.26
$ws.Range("N7").Value = 2.1
$ws.Range("P7").Value = 1.94
$ws.Range("Q7").Value = 1.73
$ws.Range("R7").Value = 1.26
$ws.Range("U7").Value = 1.01
$ws.Range("X7").Value = 20
$ws.Range("Y7").Value = 29
$ws.Range("AB7").Value = 12.5
$ws.Range("AC7").Value = 11.5
$ws.Range("AD7").Value = 32
$ws.Range("AG7").Value = 1000
$ws.Range("AK7").Value = 26
$ws.Range("AL7").Value = 50
